$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the sheet to reflect the new "through" date
$ws.Name = "Through 2022-06-16"

# Update the header label in I1 to match the new "through" date
$ws.Range("I1").Value = "2022 (through 06-16)"

# Update June (row 7) 2022 value
$ws.Range("I7").Value = 74

# Update Total (row 14) 2022 value
$ws.Range("I14").Value = 737
